$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 343; this pushes old rows 343..406 down to 344..407,
# carrying their formatting (including the date style on column D) with them.
$ws.Rows.Item(343).Insert()

# Populate the newly-inserted row 343 with the new data record.
$ws.Range("A343").Value = 9
$ws.Range("B343").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C343").Value = "Metropolitana"
$ws.Range("D343").Value = 45005
$ws.Range("E343").Value = 13
$ws.Range("F343").Value = 100112043
$ws.Range("G343").Value = "Pepino ensalada"
$ws.Range("H343").Value = "Sin especificar"
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 70
$ws.Range("K343").Value = 7000
$ws.Range("L343").Value = 8000
$ws.Range("M343").Value = 7500
$ws.Range("N343").Value = "`$/caja 60 unidades"
$ws.Range("O343").Value = "Región de Arica y Parinacota"
$ws.Range("P343").Value = 125
$ws.Range("Q343").Value = 60
$ws.Range("R343").Value = "Hortaliza"
